# Collected 7 data points
# Populate Sheet1 with a GitHub-repository data-analysis table:
# header row + 7 repositories worth of commit / contributor / release stats.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------------
$headers = @(
    "User/Repository",
    "# of Commits",
    "Avg. Time Between Commits (H)",
    "# of Contributors",
    "# of Lurkers",
    "Repo File Size (KB)",
    "# of Releases",
    "Avg. Time Between Releases (D)"
)
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# --- Data rows --------------------------------------------------------------
# User/Repository, # Commits, Avg Time Between Commits (H), # Contributors,
# # Lurkers, Repo File Size (KB), # Releases, Avg Time Between Releases (D)
$data = @(
    @("GhostPack/Seatbelt",    186,   124,                  8,   5,   940,    "N/A", "N/A"),
    @("discordjs/discord.js",  4970,  10,                   269, 263, 71800,  30,    54),
    @("davidojoy/BetterJoy",   231,   117,                  24,  23,  15204,  12,    79),
    @("citra-emu/citra",       8971,  7.4160000000000004,   198, 187, 68564,  "N/A", "N/A"),
    @("apache/druid",          10907, 6.7830000000000004,   363, 353, 226126, 30,    65),
    @("apache/dubbo",          4602,  18,                   369, 362, 30969,  30,    78.16),
    @("apache/cayenne",        667,   18.616,               33,  31,  88954,  "N/A", "N/A")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

# --- Column widths (characters, closest reachable to the authored sizes) ---
$ws.Columns.Item(1).ColumnWidth = 23.5546875
$ws.Columns.Item(2).ColumnWidth = 14.5546875
$ws.Columns.Item(3).ColumnWidth = 27.77734375
$ws.Columns.Item(4).ColumnWidth = 16.6640625
$ws.Columns.Item(5).ColumnWidth = 13
$ws.Columns.Item(6).ColumnWidth = 16.88671875
$ws.Columns.Item(7).ColumnWidth = 12.77734375
$ws.Columns.Item(8).ColumnWidth = 26.77734375
$ws.Range("I1:K1").ColumnWidth = 8.44140625

# --- Selection, matching the saved cursor position in the workbook --------
$null = $ws.Range("A13").Select()
